$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -330

$ws.Range("H98").Value = 1312.2667
$ws.Range("I98").Value = 1306.5385
$ws.Range("J98").Value = 1349.5
$ws.Range("K98").Value = 1306.5385
$ws.Range("L98").Value = 1349.5
$ws.Range("M98").Value = 191.4614999999999
$ws.Range("N98").Value = -4345.5

$ws.Range("H101").Value = 747
$ws.Range("I101").Value = 696
$ws.Range("K101").Value = 2088
$ws.Range("M101").Value = -466

$ws.Range("H103").Value = 5000
$ws.Range("J103").Value = 5000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -16172

$ws.Range("H111").Value = 1109.2142
$ws.Range("I111").Value = 849.75
$ws.Range("J111").Value = 2666
$ws.Range("K111").Value = 2549.25
$ws.Range("L111").Value = 7998
$ws.Range("M111").Value = 517.75
$ws.Range("N111").Value = -14132

$ws.Range("H113").Value = 4134.8887
$ws.Range("J113").Value = 3566.6667
$ws.Range("L113").Value = 3566.6667
$ws.Range("N113").Value = -10074.6667

$ws.Range("H115").Value = 7226.857
$ws.Range("I115").Value = 7226.857
$ws.Range("K115").Value = 21680.571
$ws.Range("M115").Value = -20113.571

$ws.Range("H116").Value = 8999.666999999999
$ws.Range("J116").Value = 8999.666999999999
$ws.Range("L116").Value = 8999.666999999999
$ws.Range("N116").Value = -15883.667

$ws.Range("H122").Value = 1312.2667
$ws.Range("I122").Value = 1306.5385
$ws.Range("J122").Value = 1349.5
$ws.Range("K122").Value = 3919.6155
$ws.Range("L122").Value = 4048.5
$ws.Range("M122").Value = -1469.6155
$ws.Range("N122").Value = -8948.5

$ws.Range("H137").Value = 4156.5415
$ws.Range("I137").Value = 1600.1
$ws.Range("J137").Value = 5982.5713
$ws.Range("K137").Value = 4800.299999999999
$ws.Range("L137").Value = 17947.7139
$ws.Range("M137").Value = -2250.299999999999
$ws.Range("N137").Value = -23047.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1279.8
$ws.Range("I2").Value = 833
$ws.Range("K2").Value = 833
$ws.Range("M2").Value = -720

$ws.Range("H61").Value = 1933.25
$ws.Range("I61").Value = 1933.25
$ws.Range("K61").Value = 1933.25
$ws.Range("M61").Value = -1721.25

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws.Range("H116").Value = 1279.8
$ws.Range("I116").Value = 833
$ws.Range("K116").Value = 833
$ws.Range("M116").Value = 1461

$ws.Range("H122").Value = 3629.7058
$ws.Range("I122").Value = 2750.3333
$ws.Range("K122").Value = 8250.999899999999
$ws.Range("M122").Value = -5800.999899999999

$ws.Range("H136").Value = 1933.25
$ws.Range("I136").Value = 1933.25
$ws.Range("K136").Value = 5799.75
$ws.Range("M136").Value = -3249.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1279.8
$ws.Range("I3").Value = 833
$ws.Range("K3").Value = 833
$ws.Range("M3").Value = -719

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4007146.2
$ws.Range("I6").Value = 4007146.2
$ws.Range("K6").Value = 4007146.2
$ws.Range("M6").Value = -4007033.2

$ws.Range("H31").Value = 4064.923
$ws.Range("I31").Value = 2793.5386
$ws.Range("K31").Value = 2793.5386
$ws.Range("M31").Value = -2498.5386

$ws.Range("H34").Value = 4064.923
$ws.Range("I34").Value = 2793.5386
$ws.Range("K34").Value = 2793.5386
$ws.Range("M34").Value = -2591.5386

$ws.Range("H99").Value = 12418.594
$ws.Range("I99").Value = 10160.3
$ws.Range("J99").Value = 16182.417
$ws.Range("K99").Value = 10160.3
$ws.Range("L99").Value = 16182.417
$ws.Range("M99").Value = -8662.299999999999
$ws.Range("N99").Value = -19178.417

$ws.Range("H122").Value = 941.6
$ws.Range("I122").Value = 902.6667
$ws.Range("K122").Value = 2708.0001
$ws.Range("M122").Value = -258.0001000000002

$ws.Range("H126").Value = 12418.594
$ws.Range("I126").Value = 10160.3
$ws.Range("J126").Value = 16182.417
$ws.Range("K126").Value = 30480.9
$ws.Range("L126").Value = 48547.251
$ws.Range("M126").Value = -28010.9
$ws.Range("N126").Value = -53487.251

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 685.4286
$ws.Range("I60").Value = 866
$ws.Range("J60").Value = 550
$ws.Range("K60").Value = 2598
$ws.Range("L60").Value = 1650
$ws.Range("M60").Value = -2347
$ws.Range("N60").Value = -2152

$ws.Range("H80").Value = 5975
$ws.Range("I80").Value = 5943.75
$ws.Range("K80").Value = 17831.25
$ws.Range("M80").Value = -16895.25

$ws.Range("H83").Value = 5975
$ws.Range("I83").Value = 5943.75
$ws.Range("K83").Value = 53493.75
$ws.Range("M83").Value = -48813.75

$ws.Range("H133").Value = 16133.25
$ws.Range("I133").Value = 4500
$ws.Range("K133").Value = 13500
$ws.Range("M133").Value = -8440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = ""

$ws.Range("H92").Value = 13000
$ws.Range("I92").Value = 25000
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 25000
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -23128
$ws.Range("N92").Value = -4744

$ws.Range("H97").Value = 1282.3077
$ws.Range("I97").Value = 1337.5454
$ws.Range("K97").Value = 1337.5454
$ws.Range("M97").Value = -841.5454

$ws.Range("H113").Value = 3830
$ws.Range("I113").Value = 2405.5
$ws.Range("J113").Value = 4399.8
$ws.Range("K113").Value = 2405.5
$ws.Range("L113").Value = 4399.8
$ws.Range("M113").Value = -235.5
$ws.Range("N113").Value = -8739.799999999999

$ws.Range("H126").Value = 5431.5557
$ws.Range("I126").Value = 4440
$ws.Range("J126").Value = 5555.5
$ws.Range("K126").Value = 13320
$ws.Range("L126").Value = 16666.5
$ws.Range("M126").Value = -10850
$ws.Range("N126").Value = -21606.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4422.706
$ws.Range("I40").Value = 4327.643
$ws.Range("J40").Value = 4866.3335
$ws.Range("K40").Value = 4327.643
$ws.Range("L40").Value = 4866.3335
$ws.Range("M40").Value = -4191.643
$ws.Range("N40").Value = -5138.3335

$ws.Range("H93").Value = 1111
$ws.Range("I93").Value = 796.1667
$ws.Range("K93").Value = 796.1667
$ws.Range("M93").Value = 451.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 36007.5
$ws.Range("I51").Value = 22000
$ws.Range("K51").Value = 22000
$ws.Range("M51").Value = -21490

$ws.Range("H56").Value = 9995
$ws.Range("I56").Value = 9995
$ws.Range("K56").Value = 9995
$ws.Range("M56").Value = -9281

$ws.Range("H60").Value = 141999.4
$ws.Range("J60").Value = 149999.25
$ws.Range("L60").Value = 149999.25
$ws.Range("N60").Value = -151643.25

$ws.Range("H107").Value = 805.0909
$ws.Range("I107").Value = 642.75
$ws.Range("J107").Value = 999.9
$ws.Range("K107").Value = 1928.25
$ws.Range("L107").Value = 2999.7
$ws.Range("M107").Value = -8.25
$ws.Range("N107").Value = -6839.7
